$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume(1h) percentage values
# as scraped on Mon May 15 20:30:18 UTC 2023
$ws.Range('D2').Value = '27.543.42'
$ws.Range('D3').Value = '1.838.47'
$ws.Range('E3').Value = '  -0.60%  '
$ws.Range('E4').Value = '  -2.57%  '
$ws.Range('D5').Value = '316.56'
$ws.Range('E5').Value = '  -1.54%  '
$ws.Range('D6').Value = '1.007'
$ws.Range('E6').Value = '  -2.29%  '
$ws.Range('D7').Value = '0.4302'
$ws.Range('E7').Value = '  -2.03%  '
$ws.Range('D8').Value = '0.3717'
$ws.Range('E8').Value = '  -1.93%  '
$ws.Range('D9').Value = '0.07288'
$ws.Range('E9').Value = '  -1.40%  '
$ws.Range('D10').Value = '0.8695'
$ws.Range('E10').Value = '  -1.44%  '
$ws.Range('D11').Value = '21.24'
$ws.Range('E11').Value = '  -1.26%  '
$ws.Range('D12').Value = '1.847.85'
$ws.Range('E12').Value = '  -0.29%  '
$ws.Range('D13').Value = '6.713'
$ws.Range('E13').Value = '  +0.36%  '
$ws.Range('D14').Value = '5.374'
$ws.Range('E14').Value = '  -2.25%  '
$ws.Range('D15').Value = '0.07100'
$ws.Range('E15').Value = '  -0.89%  '
$ws.Range('E16').Value = '  +4.37%  '
$ws.Range('D17').Value = '1.010'
$ws.Range('E17').Value = '  -2.69%  '
$ws.Range('D18').Value = '0.000008957'
$ws.Range('E18').Value = '  -1.32%  '
$ws.Range('D19').Value = '1.006'
$ws.Range('E19').Value = '  -2.40%  '
$ws.Range('D20').Value = '15.31'
$ws.Range('E20').Value = '  -0.95%  '
$ws.Range('D21').Value = '27.552.34'
$ws.Range('E21').Value = '  -0.59%  '
$ws.Range('D22').Value = '5.178'
$ws.Range('E22').Value = '  -1.98%  '
$ws.Range('D23').Value = '10.97'
$ws.Range('E23').Value = '  -2.40%  '
$ws.Range('D24').Value = '2.068.91'
$ws.Range('E24').Value = '  -1.35%  '
$ws.Range('D25').Value = '2.010'
$ws.Range('E25').Value = '  -2.00%  '
$ws.Range('D26').Value = '154.10'
$ws.Range('E26').Value = '  -2.94%  '
$ws.Range('D27').Value = '18.51'
$ws.Range('E27').Value = '  -0.92%  '
$ws.Range('D28').Value = '2.161'
$ws.Range('E28').Value = '  +8.64%  '
$ws.Range('D29').Value = '5.305'
$ws.Range('E29').Value = '  -0.42%  '
$ws.Range('D30').Value = '117.46'
$ws.Range('E30').Value = '  -0.24%  '
$ws.Range('D31').Value = '0.08869'
$ws.Range('E31').Value = '  -2.27%  '
$ws.Range('D32').Value = '1.212'
$ws.Range('E32').Value = '  +0.41%  '
$ws.Range('D33').Value = '0.7705'
$ws.Range('E33').Value = '  -0.18%  '
$ws.Range('D34').Value = '4.503'
$ws.Range('E34').Value = '  -1.12%  '
$ws.Range('D35').Value = '2.907'
$ws.Range('E35').Value = '  -3.18%  '
$ws.Range('D36').Value = '1.008'
$ws.Range('E36').Value = '  -2.46%  '
$ws.Range('D37').Value = '1.127'
$ws.Range('E37').Value = '  -2.13%  '
$ws.Range('D38').Value = '0.05303'
$ws.Range('E38').Value = '  +0.83%  '
$ws.Range('D39').Value = '0.01965'
$ws.Range('E39').Value = '  -0.30%  '
$ws.Range('D40').Value = '7.162'
$ws.Range('E40').Value = '  +4.36%  '
$ws.Range('D41').Value = '2.882'
$ws.Range('E41').Value = '  +1.12%  '
$ws.Range('D42').Value = '0.5108'
$ws.Range('E42').Value = '  -1.28%  '
$ws.Range('E43').Value = '  +0.64%  '
$ws.Range('D44').Value = '8.713'
$ws.Range('E44').Value = '  +0.19%  '
$ws.Range('D45').Value = '10.57'
$ws.Range('E45').Value = '  -1.10%  '
$ws.Range('D46').Value = '106.88'
$ws.Range('E46').Value = '  -3.00%  '
$ws.Range('D47').Value = '0.4731'
$ws.Range('E47').Value = '  +0.84%  '
$ws.Range('D48').Value = '0.06436'
$ws.Range('E48').Value = '  -2.12%  '
$ws.Range('D49').Value = '1.007'
$ws.Range('E49').Value = '  -2.57%  '
$ws.Range('E50').Value = '  -1.25%  '
$ws.Range('D51').Value = '1.836'
$ws.Range('E51').Value = '  -2.50%  '
